$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 573.36365
$ws.Range("I28").Value = 430.7
$ws.Range("K28").Value = 430.7
$ws.Range("M28").Value = 54.30000000000001
$ws.Range("H70").Value = 2137.0833
$ws.Range("I70").Value = 2049.7144
$ws.Range("J70").Value = 2259.4
$ws.Range("K70").Value = 6149.1432
$ws.Range("L70").Value = 6778.200000000001
$ws.Range("M70").Value = -5879.1432
$ws.Range("N70").Value = -7318.200000000001
$ws.Range("H73").Value = 2137.0833
$ws.Range("I73").Value = 2049.7144
$ws.Range("J73").Value = 2259.4
$ws.Range("K73").Value = 6149.1432
$ws.Range("L73").Value = 6778.200000000001
$ws.Range("M73").Value = -5213.1432
$ws.Range("N73").Value = -8650.200000000001
$ws.Range("H82").Value = 514.5
$ws.Range("I82").Value = 514.5
$ws.Range("K82").Value = 1543.5
$ws.Range("M82").Value = -1137.5
$ws.Range("H85").Value = 514.5
$ws.Range("I85").Value = 514.5
$ws.Range("K85").Value = 1543.5
$ws.Range("M85").Value = -139.5
$ws.Range("H111").Value = 3299.5
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("H112").Value = 2495.8
$ws.Range("J112").Value = 2920
$ws.Range("L112").Value = 8760
$ws.Range("N112").Value = -10976
$ws.Range("H118").Value = 457.6
$ws.Range("I118").Value = 457.6
$ws.Range("K118").Value = 1372.8
$ws.Range("M118").Value = 284.1999999999998
$ws.Range("H138").Value = 4354.154
$ws.Range("J138").Value = 3691.5454
$ws.Range("L138").Value = 11074.6362
$ws.Range("N138").Value = -21354.6362
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 12016.333
$ws.Range("J34").Value = 29999
$ws.Range("L34").Value = 29999
$ws.Range("N34").Value = -30541
$ws.Range("H45").Value = 6600
$ws.Range("I45").Value = 6000
$ws.Range("K45").Value = 6000
$ws.Range("M45").Value = -5623
$ws.Range("H110").Value = 2096.8
$ws.Range("I110").Value = 2096.8
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 2096.8
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -51.80000000000018
$ws.Range("N110").ClearContents()
$ws.Range("H122").Value = 2425
$ws.Range("I122").Value = 2750
$ws.Range("K122").Value = 8250
$ws.Range("M122").Value = -5800
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 63780
$ws.Range("J58").Value = 63780
$ws.Range("L58").Value = 63780
$ws.Range("N58").Value = -64368
$ws.Range("H74").Value = 55387.5
$ws.Range("J74").Value = 55387.5
$ws.Range("L74").Value = 55387.5
$ws.Range("N74").Value = -57259.5
$ws.Range("H77").Value = 55387.5
$ws.Range("J77").Value = 55387.5
$ws.Range("L77").Value = 166162.5
$ws.Range("N77").Value = -175522.5
$ws.Range("H86").Value = 3998.3333
$ws.Range("I86").Value = 3997.5
$ws.Range("K86").Value = 3997.5
$ws.Range("M86").Value = -2874.5
$ws.Range("H89").Value = 3998.3333
$ws.Range("I89").Value = 3997.5
$ws.Range("K89").Value = 19987.5
$ws.Range("M89").Value = -14371.5
$ws.Range("H139").Value = 74999
$ws.Range("J139").Value = 74999
$ws.Range("L139").Value = 74999
$ws.Range("N139").Value = -85279
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 566.3333
$ws.Range("I3").Value = 566.3333
$ws.Range("K3").Value = 566.3333
$ws.Range("M3").Value = -453.3333
$ws.Range("H18").Value = 79955
$ws.Range("J18").Value = 79955
$ws.Range("L18").Value = 79955
$ws.Range("N18").Value = -80415
$ws.Range("H22").Value = 1640.7273
$ws.Range("I22").Value = 1024.6666
$ws.Range("K22").Value = 1024.6666
$ws.Range("M22").Value = -674.6666
$ws.Range("H43").Value = 22833.334
$ws.Range("J43").Value = 22833.334
$ws.Range("L43").Value = 22833.334
$ws.Range("N43").Value = -23201.334
$ws.Range("H63").Value = 96665.664
$ws.Range("J63").Value = 96665.664
$ws.Range("L63").Value = 96665.664
$ws.Range("N63").Value = -98037.664
$ws.Range("H66").Value = 96665.664
$ws.Range("J66").Value = 96665.664
$ws.Range("L66").Value = 289996.992
$ws.Range("N66").Value = -296860.992
$ws.Range("H101").Value = 22833.334
$ws.Range("J101").Value = 22833.334
$ws.Range("L101").Value = 22833.334
$ws.Range("N101").Value = -29323.334
$ws.Range("H132").Value = 6645.609
$ws.Range("I132").Value = 5892.45
$ws.Range("J132").Value = 11666.667
$ws.Range("K132").Value = 17677.35
$ws.Range("L132").Value = 35000.001
$ws.Range("M132").Value = -15147.35
$ws.Range("N132").Value = -40060.001
$ws.Range("H134").Value = 2722
$ws.Range("I134").Value = 2999.5
$ws.Range("J134").Value = 2167
$ws.Range("K134").Value = 8998.5
$ws.Range("L134").Value = 6501
$ws.Range("M134").Value = -6463.5
$ws.Range("N134").Value = -11571
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 327.33334
$ws.Range("I8").Value = 327.33334
$ws.Range("K8").Value = 982.0000200000001
$ws.Range("M8").Value = -843.0000200000001
$ws.Range("H34").Value = 665
$ws.Range("J34").Value = 1026
$ws.Range("L34").Value = 3078
$ws.Range("N34").Value = -3246
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H55").Value = 1070.6364
$ws.Range("J55").Value = 1100
$ws.Range("L55").Value = 3300
$ws.Range("N55").Value = -3654
$ws.Range("H122").Value = 1977
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("H102").Value = 1691.8125
$ws.Range("I102").Value = 961.7273
$ws.Range("K102").Value = 961.7273
$ws.Range("M102").Value = 660.2727
$ws.Range("H107").Value = 773.875
$ws.Range("I107").Value = 415.16666
$ws.Range("K107").Value = 415.16666
$ws.Range("M107").Value = 1504.83334
$ws.Range("H113").Value = 1000
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2400
$ws.Range("I46").Value = 2362.75
$ws.Range("J46").Value = 2421.2856
$ws.Range("K46").Value = 2362.75
$ws.Range("L46").Value = 2421.2856
$ws.Range("M46").Value = -2174.75
$ws.Range("N46").Value = -2797.2856
$ws.Range("H61").Value = 3600
$ws.Range("I61").Value = 4333.3335
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 4333.3335
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -4131.3335
$ws.Range("N61").Value = -2904
$ws.Range("H93").Value = 1000
$ws.Range("I93").Value = 1000
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1000
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 248
$ws.Range("N93").ClearContents()
$ws.Range("H113").Value = 3600
$ws.Range("I113").Value = 4333.3335
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 4333.3335
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = -2163.3335
$ws.Range("N113").Value = -6840
$ws.Range("H122").Value = 6997
$ws.Range("I122").Value = 6332.6665
$ws.Range("J122").Value = 8990
$ws.Range("K122").Value = 18997.9995
$ws.Range("L122").Value = 26970
$ws.Range("M122").Value = -16547.9995
$ws.Range("N122").Value = -31870
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7999.75
$ws.Range("I62").Value = 7999.75
$ws.Range("K62").Value = 7999.75
$ws.Range("M62").Value = -7375.75
$ws.Range("H65").Value = 7999.75
$ws.Range("I65").Value = 7999.75
$ws.Range("K65").Value = 39998.75
$ws.Range("M65").Value = -36878.75
$ws.Range("H81").Value = 2800.2
$ws.Range("I81").Value = 1501
$ws.Range("J81").Value = 3666.3333
$ws.Range("K81").Value = 3002
$ws.Range("L81").Value = 7332.6666
$ws.Range("M81").Value = -1941
$ws.Range("N81").Value = -9454.6666
$ws.Range("H84").Value = 2800.2
$ws.Range("I84").Value = 1501
$ws.Range("J84").Value = 3666.3333
$ws.Range("K84").Value = 15010
$ws.Range("L84").Value = 36663.333
$ws.Range("M84").Value = -9706
$ws.Range("N84").Value = -47271.333
$ws.Range("H122").Value = 200
$ws.Range("I122").Value = 200
$ws.Range("K122").Value = 600
$ws.Range("M122").Value = 1850
$ws.Range("H132").Value = 201169.2
$ws.Range("J132").Value = 1000
$ws.Range("L132").Value = 3000
$ws.Range("N132").Value = -8060
